$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto list data (coin order, prices and 1h volume changes refreshed)

# Row 2
$ws.Cells.Item(2, 5).Value = '  -6.04%  '
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '62.191.78'
$ws.Cells.Item(2, 4).Style = "Normal"

# Row 3
$ws.Cells.Item(3, 5).Value = '  -6.58%  '
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '3.176.36'
$ws.Cells.Item(3, 4).Style = "Normal"

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.08%  '

# Row 5
$ws.Cells.Item(5, 2).Value = 'BNB'
$ws.Cells.Item(5, 3).Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Cells.Item(5, 5).Value = '  -4.55%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '509.70'
$ws.Cells.Item(5, 4).Style = "Normal"

# Row 6
$ws.Cells.Item(6, 2).Value = 'Solana'
$ws.Cells.Item(6, 3).Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Cells.Item(6, 5).Value = '  -8.88%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '170.63'
$ws.Cells.Item(6, 4).Style = "Normal"

# Row 7
$ws.Cells.Item(7, 5).Value = '  -5.99%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.582'
$ws.Cells.Item(7, 4).Style = "Normal"

# Row 8
$ws.Cells.Item(8, 5).Value = '  +0.03%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  -6.48%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '3.178.34'
$ws.Cells.Item(9, 4).Style = "Normal"

# Row 10
$ws.Cells.Item(10, 5).Value = '  -6.51%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.587'
$ws.Cells.Item(10, 4).Style = "Normal"

# Row 11
$ws.Cells.Item(11, 5).Value = '  -9.83%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '52.13'
$ws.Cells.Item(11, 4).Style = "Normal"

# Row 12
$ws.Cells.Item(12, 5).Value = '  -6.27%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.128'
$ws.Cells.Item(12, 4).Style = "Normal"

# Row 13
$ws.Cells.Item(13, 5).Value = '  -2.89%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.0000247'
$ws.Cells.Item(13, 4).Style = "Normal"

# Row 14
$ws.Cells.Item(14, 5).Value = '  -6.57%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '8.74'
$ws.Cells.Item(14, 4).Style = "Normal"

# Row 15
$ws.Cells.Item(15, 5).Value = '  -6.52%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '3.700.56'
$ws.Cells.Item(15, 4).Style = "Normal"

# Row 16
$ws.Cells.Item(16, 2).Value = 'WrappedEther'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(16, 5).Value = '  -7.14%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '3.181.12'
$ws.Cells.Item(16, 4).Style = "Normal"

# Row 17
$ws.Cells.Item(17, 2).Value = 'TRON'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Cells.Item(17, 5).Value = '  -8.31%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '0.113'
$ws.Cells.Item(17, 4).Style = "Normal"

# Row 18
$ws.Cells.Item(18, 5).Value = '  -5.87%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '62.199.89'
$ws.Cells.Item(18, 4).Style = "Normal"

# Row 19
$ws.Cells.Item(19, 5).Value = '  -3.48%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '16.87'
$ws.Cells.Item(19, 4).Style = "Normal"

# Row 20
$ws.Cells.Item(20, 5).Value = '  -5.20%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '10.72'
$ws.Cells.Item(20, 4).Style = "Normal"

# Row 21
$ws.Cells.Item(21, 5).Value = '  -5.07%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '0.934'
$ws.Cells.Item(21, 4).Style = "Normal"

# Row 22
$ws.Cells.Item(22, 5).Value = '  -6.20%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '358.82'
$ws.Cells.Item(22, 4).Style = "Normal"

# Row 23
$ws.Cells.Item(23, 5).Value = '  -2.09%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '10.94'
$ws.Cells.Item(23, 4).Style = "Normal"

# Row 24
$ws.Cells.Item(24, 5).Value = '  -3.48%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '3.64'
$ws.Cells.Item(24, 4).Style = "Normal"

# Row 25
$ws.Cells.Item(25, 5).Value = '  -5.00%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '78.84'
$ws.Cells.Item(25, 4).Style = "Normal"

# Row 26
$ws.Cells.Item(26, 5).Value = '  +1.93%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '6.18'
$ws.Cells.Item(26, 4).Style = "Normal"

# Row 27
$ws.Cells.Item(27, 5).Value = '  +1.83%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '3.73'
$ws.Cells.Item(27, 4).Style = "Normal"

# Row 28
$ws.Cells.Item(28, 5).Value = '  -4.76%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '2.57'
$ws.Cells.Item(28, 4).Style = "Normal"

# Row 29
$ws.Cells.Item(29, 5).Value = '  -4.34%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '10.99'
$ws.Cells.Item(29, 4).Style = "Normal"

# Row 30
$ws.Cells.Item(30, 5).Value = '  -5.57%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '7.98'
$ws.Cells.Item(30, 4).Style = "Normal"

# Row 31
$ws.Cells.Item(31, 2).Value = 'Bittensor'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(31, 5).Value = '  -7.04%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '638.05'
$ws.Cells.Item(31, 4).Style = "Normal"

# Row 32
$ws.Cells.Item(32, 2).Value = 'EthereumClassic'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(32, 5).Value = '  -6.39%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '27.73'
$ws.Cells.Item(32, 4).Style = "Normal"

# Row 33
$ws.Cells.Item(33, 5).Value = '  -6.14%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '6.19'
$ws.Cells.Item(33, 4).Style = "Normal"

# Row 34
$ws.Cells.Item(34, 5).Value = '  -2.85%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '10.89'
$ws.Cells.Item(34, 4).Style = "Normal"

# Row 35
$ws.Cells.Item(35, 5).Value = '  -4.04%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.102'
$ws.Cells.Item(35, 4).Style = "Normal"

# Row 36
$ws.Cells.Item(36, 5).Value = '  -8.00%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '57.00'
$ws.Cells.Item(36, 4).Style = "Normal"

# Row 37
$ws.Cells.Item(37, 5).Value = '  +0.26%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  -2.15%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '35.76'
$ws.Cells.Item(38, 4).Style = "Normal"

# Row 39
$ws.Cells.Item(39, 5).Value = '  -1.72%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.366'
$ws.Cells.Item(39, 4).Style = "Normal"

# Row 40
$ws.Cells.Item(40, 5).Value = '  +0.16%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  +8.24%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.0₃0674'
$ws.Cells.Item(41, 4).Style = "Normal"

# Row 42
$ws.Cells.Item(42, 2).Value = 'Maker'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(42, 5).Value = '  -3.13%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '2.829.41'
$ws.Cells.Item(42, 4).Style = "Normal"

# Row 43
$ws.Cells.Item(43, 2).Value = 'Kaspa'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(43, 5).Value = '  -4.90%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.120'
$ws.Cells.Item(43, 4).Style = "Normal"

# Row 44
$ws.Cells.Item(44, 5).Value = '  +3.58%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '2.48'
$ws.Cells.Item(44, 4).Style = "Normal"

# Row 45
$ws.Cells.Item(45, 5).Value = '  +0.36%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '2.63'
$ws.Cells.Item(45, 4).Style = "Normal"

# Row 46
$ws.Cells.Item(46, 5).Value = '  -10.21%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '2.53'
$ws.Cells.Item(46, 4).Style = "Normal"

# Row 47
$ws.Cells.Item(47, 5).Value = '  -1.93%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '0.0378'
$ws.Cells.Item(47, 4).Style = "Normal"

# Row 48
$ws.Cells.Item(48, 2).Value = 'Monero'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(48, 5).Value = '  -1.55%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '134.64'
$ws.Cells.Item(48, 4).Style = "Normal"

# Row 49
$ws.Cells.Item(49, 2).Value = 'Stacks'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(49, 5).Value = '  +4.89%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '2.69'
$ws.Cells.Item(49, 4).Style = "Normal"

# Row 50
$ws.Cells.Item(50, 2).Value = 'ApeXProtocol'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Cells.Item(50, 5).Value = '  -1.16%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '2.87'
$ws.Cells.Item(50, 4).Style = "Normal"

# Row 51
$ws.Cells.Item(51, 5).Value = '  -4.75%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.121'
$ws.Cells.Item(51, 4).Style = "Normal"
